$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32 corresponds to the match: Sporting Cristal 0-1 U. de Deportes (24/10/2025)
$ws.Range("A32").Value = "24/10/2025"
$ws.Range("B32").Value = "Sporting Cristal"
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 1
$ws.Range("E32").Value = "U. de Deportes"
$ws.Range("F32").Value = "L"
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1
$ws.Range("K32").Value = 0.8
$ws.Range("L32").Value = 1.32
$ws.Range("M32").Value = 16
$ws.Range("N32").Value = 8
$ws.Range("O32").Value = 4
$ws.Range("P32").Value = 3
